$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 31615
$ws.Range("I32").Value = 106248.5
$ws.Range("J32").Value = 4475.5454
$ws.Range("K32").Value = 106248.5
$ws.Range("L32").Value = 4475.5454
$ws.Range("M32").Value = -105922.5
$ws.Range("N32").Value = -5127.5454
$ws.Range("H41").Value = 125127.75
$ws.Range("I41").Value = 125
$ws.Range("J41").Value = 166795.33
$ws.Range("K41").Value = 125
$ws.Range("L41").Value = 166795.33
$ws.Range("M41").Value = 315
$ws.Range("N41").Value = -167675.33
$ws.Range("H53").Value = 4360.4
$ws.Range("I53").Value = 244.75
$ws.Range("J53").Value = 9064
$ws.Range("K53").Value = 244.75
$ws.Range("L53").Value = 9064
$ws.Range("M53").Value = 392.25
$ws.Range("N53").Value = -10338
$ws.Range("H112").Value = 30477.914
$ws.Range("J112").Value = 1966.6923
$ws.Range("L112").Value = 5900.0769
$ws.Range("N112").Value = -8116.0769
$ws.Range("H116").Value = 32469908
$ws.Range("I116").Value = 35859500
$ws.Range("K116").Value = 35859500
$ws.Range("M116").Value = -35856058
$ws.Range("H132").Value = 2570.875
$ws.Range("I132").Value = 2473.9167
$ws.Range("J132").Value = 3443.5
$ws.Range("K132").Value = 7421.750100000001
$ws.Range("L132").Value = 10330.5
$ws.Range("M132").Value = -4891.750100000001
$ws.Range("N132").Value = -15390.5
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 100000
$ws.Range("N134").Value = -110140
$ws.Range("H137").Value = 71429600
$ws.Range("I137").Value = 83334370
$ws.Range("K137").Value = 250003110
$ws.Range("M137").Value = -250000560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1249.5
$ws.Range("I2").Value = 999.3333
$ws.Range("K2").Value = 999.3333
$ws.Range("M2").Value = -886.3333
$ws.Range("H61").Value = 13890273
$ws.Range("I61").Value = 13890273
$ws.Range("K61").Value = 13890273
$ws.Range("M61").Value = -13890061
$ws.Range("H110").Value = 83418430
$ws.Range("I110").Value = 125063910
$ws.Range("K110").Value = 125063910
$ws.Range("M110").Value = -125061865
$ws.Range("H116").Value = 1249.5
$ws.Range("I116").Value = 999.3333
$ws.Range("K116").Value = 999.3333
$ws.Range("M116").Value = 1294.6667
$ws.Range("H132").Value = 25001776
$ws.Range("J132").Value = 3700
$ws.Range("L132").Value = 11100
$ws.Range("N132").Value = -16160
$ws.Range("H136").Value = 13890273
$ws.Range("I136").Value = 13890273
$ws.Range("K136").Value = 41670819
$ws.Range("M136").Value = -41668269

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1249.5
$ws.Range("I3").Value = 999.3333
$ws.Range("K3").Value = 999.3333
$ws.Range("M3").Value = -885.3333
$ws.Range("H20").Value = 18562.818
$ws.Range("I20").Value = 24172.732
$ws.Range("K20").Value = 24172.732
$ws.Range("M20").Value = -23925.732
$ws.Range("H22").Value = 189.14285
$ws.Range("I22").Value = 137.33333
$ws.Range("K22").Value = 137.33333
$ws.Range("M22").Value = 35.66667000000001
$ws.Range("H75").Value = 19142.2
$ws.Range("I75").Value = 18927.75
$ws.Range("K75").Value = 18927.75
$ws.Range("M75").Value = -17991.75
$ws.Range("H78").Value = 19142.2
$ws.Range("I78").Value = 18927.75
$ws.Range("K78").Value = 56783.25
$ws.Range("M78").Value = -52103.25
$ws.Range("H96").Value = 24231.75
$ws.Range("I96").Value = 24231.75
$ws.Range("K96").Value = 24231.75
$ws.Range("M96").Value = -21485.75
$ws.Range("H97").Value = 7831.4
$ws.Range("I97").Value = 7831.4
$ws.Range("K97").Value = 7831.4
$ws.Range("M97").Value = -6840.4
$ws.Range("H105").Value = 1957.8125
$ws.Range("I105").Value = 1892.9166
$ws.Range("K105").Value = 1892.9166
$ws.Range("M105").Value = -145.9166
$ws.Range("H107").Value = 21752844
$ws.Range("I107").Value = 14718.571
$ws.Range("J107").Value = 250003150
$ws.Range("K107").Value = 14718.571
$ws.Range("L107").Value = 250003150
$ws.Range("M107").Value = -12798.571
$ws.Range("N107").Value = -250006990

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1010.1667
$ws.Range("I16").Value = 1010.1667
$ws.Range("K16").Value = 1010.1667
$ws.Range("M16").Value = -723.1667
$ws.Range("H107").Value = 2585.9473
$ws.Range("I107").Value = 2487.5
$ws.Range("J107").Value = 2861.6
$ws.Range("K107").Value = 2487.5
$ws.Range("L107").Value = 2861.6
$ws.Range("M107").Value = -567.5
$ws.Range("N107").Value = -6701.6
$ws.Range("H113").Value = 1010.1667
$ws.Range("I113").Value = 1010.1667
$ws.Range("K113").Value = 1010.1667
$ws.Range("M113").Value = 1159.8333
$ws.Range("H122").Value = 2208.6365
$ws.Range("I122").Value = 1816.8
$ws.Range("J122").Value = 3048.2856
$ws.Range("K122").Value = 5450.4
$ws.Range("L122").Value = 9144.856800000001
$ws.Range("M122").Value = -3000.4
$ws.Range("N122").Value = -14044.8568
$ws.Range("H132").Value = 3293.0588
$ws.Range("I132").Value = 2992.7856
$ws.Range("K132").Value = 8978.356800000001
$ws.Range("M132").Value = -6448.356800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 187574.44
$ws.Range("I11").Value = 54921.58
$ws.Range("K11").Value = 164764.74
$ws.Range("M11").Value = -164624.74
$ws.Range("H17").Value = 486.14285
$ws.Range("I17").Value = 437.66666
$ws.Range("K17").Value = 1312.99998
$ws.Range("M17").Value = -1143.99998
$ws.Range("H113").Value = 1045.7778
$ws.Range("J113").Value = 1175.5714
$ws.Range("L113").Value = 3526.7142
$ws.Range("N113").Value = -7866.7142
$ws.Range("H122").Value = 399.1905
$ws.Range("I122").Value = 325.54544
$ws.Range("J122").Value = 480.2
$ws.Range("K122").Value = 2929.90896
$ws.Range("L122").Value = 4321.8
$ws.Range("M122").Value = -479.9089599999998
$ws.Range("N122").Value = -9221.799999999999
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 24803.5
$ws.Range("J21").Value = 24803.5
$ws.Range("L21").Value = 24803.5
$ws.Range("N21").Value = -25149.5
$ws.Range("H30").Value = 24803.5
$ws.Range("J30").Value = 24803.5
$ws.Range("L30").Value = 24803.5
$ws.Range("N30").Value = -25013.5
$ws.Range("H113").Value = 1388
$ws.Range("I113").Value = 986.9231
$ws.Range("J113").Value = 3995
$ws.Range("K113").Value = 986.9231
$ws.Range("L113").Value = 3995
$ws.Range("M113").Value = 1183.0769
$ws.Range("N113").Value = -8335
$ws.Range("H122").Value = 2582.4211
$ws.Range("I122").Value = 2148.6428
$ws.Range("J122").Value = 3797
$ws.Range("K122").Value = 6445.928400000001
$ws.Range("L122").Value = 11391
$ws.Range("M122").Value = -3995.928400000001
$ws.Range("N122").Value = -16291
$ws.Range("H126").Value = 10200.4375
$ws.Range("I126").Value = 13091.546
$ws.Range("K126").Value = 39274.638
$ws.Range("M126").Value = -36804.638
$ws.Range("H132").Value = 3152.457
$ws.Range("I132").Value = 2154.9092
$ws.Range("K132").Value = 6464.7276
$ws.Range("M132").Value = -3934.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 83334950
$ws.Range("I7").Value = 83334950
$ws.Range("K7").Value = 83334950
$ws.Range("M7").Value = -83334838
$ws.Range("H39").Value = 22500
$ws.Range("I39").Value = 20000
$ws.Range("J39").Value = 25000
$ws.Range("K39").Value = 20000
$ws.Range("L39").Value = 25000
$ws.Range("M39").Value = -19540
$ws.Range("N39").Value = -25920
$ws.Range("H46").Value = 1774.5238
$ws.Range("I46").Value = 893
$ws.Range("J46").Value = 2435.6667
$ws.Range("K46").Value = 893
$ws.Range("L46").Value = 2435.6667
$ws.Range("M46").Value = -705
$ws.Range("N46").Value = -2811.6667
$ws.Range("H61").Value = 22239.584
$ws.Range("J61").Value = 21540.777
$ws.Range("L61").Value = 21540.777
$ws.Range("N61").Value = -21944.777
$ws.Range("H113").Value = 22239.584
$ws.Range("J113").Value = 21540.777
$ws.Range("L113").Value = 21540.777
$ws.Range("N113").Value = -25880.777
$ws.Range("H126").Value = 83334950
$ws.Range("I126").Value = 83334950
$ws.Range("K126").Value = 250004850
$ws.Range("M126").Value = -250002380

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6468.8667
$ws.Range("I81").Value = 2452.889
$ws.Range("K81").Value = 4905.778
$ws.Range("M81").Value = -3844.778
$ws.Range("H84").Value = 6468.8667
$ws.Range("I84").Value = 2452.889
$ws.Range("K84").Value = 24528.89
$ws.Range("M84").Value = -19224.89
$ws.Range("H100").Value = 3359.2222
$ws.Range("I100").Value = 4097.5713
$ws.Range("K100").Value = 8195.142599999999
$ws.Range("M100").Value = -7654.142599999999
$ws.Range("H107").Value = 483.76923
$ws.Range("I107").Value = 417.18182
$ws.Range("J107").Value = 850
$ws.Range("K107").Value = 1251.54546
$ws.Range("L107").Value = 2550
$ws.Range("M107").Value = 668.45454
$ws.Range("N107").Value = -6390
$ws.Range("H113").Value = 328.2
$ws.Range("I113").Value = 226.91667
$ws.Range("K113").Value = 680.75001
$ws.Range("M113").Value = 1489.24999
$ws.Range("H136").Value = 6109.8
$ws.Range("I136").Value = 2516.5
$ws.Range("K136").Value = 7549.5
$ws.Range("M136").Value = -4999.5
